# SimulatorScenarios.xlsx: bump agent_num (column D) for both scenario rows
# from 1000 to 300, then leave the selection on D4 (matching the saved
# cursor position in the edited workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("simulator_scenarios")
$ws.Activate()

$ws.Range("D2").Value = 300
$ws.Range("D3").Value = 300

$ws.Range("D4").Select()
